$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("chart")
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$sc = $chart.SeriesCollection()
$news = $sc.NewSeries()
$news.Formula = "=SERIES(data!`$A`$16,data!`$AS`$2:`$ZZ`$2,data!`$AS`$16:`$ZZ`$16,13)"
$news.MarkerStyle = -4142
$news.Smooth = $false
try { $news.Border.Weight = 2.25; Write-Host "borderweight ok" } catch { Write-Host "ERR bw: $_" }
try { $news.Border.Color = 255; Write-Host "bordercolor ok" } catch { Write-Host "ERR bc: $_" }
